$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (kept as literal text, matching the
# workbook's inline-string cells for Price / Volume(1h) columns).
$updates = [ordered]@{
    "D2" = "321.04"
    "E2" = "-3.22%"
    "D3" = "42.76"
    "E3" = "-6.12%"
    "D4" = "5.202"
    "D5" = "0.08169"
    "E5" = "-2.31%"
    "D6" = "4.316"
    "E6" = "-3.01%"
    "E7" = "-13.40%"
    "D8" = "0.9338"
    "E8" = "-5.70%"
    "D9" = "0.1110"
    "E9" = "-7.27%"
    "D10" = "0.1854"
    "E10" = "-3.89%"
    "D11" = "0.04713"
    "E11" = "0.96%"
    "D12" = "0.09363"
    "E12" = "-5.37%"
    "D13" = "7.477"
    "E13" = "-27.47%"
    "E14" = "-0.24%"
    "D15" = "0.001293"
    "E15" = "-0.32%"
    "D16" = "0.005940"
    "E16" = "0.28%"
    "D17" = "3.355"
    "E19" = "-0.60%"
    "E20" = "0.64%"
    "E21" = "-0.59%"
    "D22" = "0.04149"
    "E22" = "-0.08%"
    "D23" = "0.001246"
    "E23" = "-3.61%"
    "D24" = "0.004300"
    "E24" = "-5.15%"
    "D25" = "0.0001201"
    "D38" = "0.02702"
    "E38" = "0.09%"
    "D39" = "0.05542"
    "E39" = "-3.51%"
    "D40" = "0.008045"
    "E40" = "2.02%"
    "D41" = "0.1396"
    "E41" = "-2.59%"
    "D42" = "0.006549"
    "E42" = "-16.78%"
    "D43" = "0.002087"
    "E43" = "3.25%"
    "D44" = "0.008263"
    "E44" = "-7.47%"
    "D45" = "0.3494"
    "E45" = "2.55%"
    "D46" = "0.00006920"
    "E46" = "-1.64%"
    "E47" = "-0.01%"
    "D48" = "0.003340"
    "E48" = "-1.91%"
    "E49" = "0.03%"
    "E50" = "-0.01%"
    "E51" = "-0.01%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force a Text number format first so Excel does not reinterpret
    # numeric- or percent-looking strings as numbers, then restore the
    # cell's original (default) style so no formatting is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
